# Applies the scheduled-runner profit recalculation to the Leve profit sheets.
# Source values come from an updated Market Board price pull; formulas are untouched,
# only the cached price/profit cell values are refreshed per row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5825
$ws.Range("I113").Value = 6981.1665
$ws.Range("K113").Value = 6981.1665
$ws.Range("M113").Value = -3727.1665
$ws.Range("H121").Value = 4315.353
$ws.Range("J121").Value = 4430.125
$ws.Range("L121").Value = 13290.375
$ws.Range("N121").Value = -16784.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 3000
$ws.Range("J15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("N15").Value = -3700
$ws.Range("H55").Value = 30587.908
$ws.Range("J55").Value = 30587.908
$ws.Range("L55").Value = 30587.908
$ws.Range("N55").Value = -31217.908
$ws.Range("H80").Value = 38876
$ws.Range("J80").Value = 37306.668
$ws.Range("L80").Value = 37306.668
$ws.Range("N80").Value = -39302.668
$ws.Range("H83").Value = 38876
$ws.Range("J83").Value = 37306.668
$ws.Range("L83").Value = 111920.004
$ws.Range("N83").Value = -121904.004
$ws.Range("H122").Value = 2588.258
$ws.Range("I122").Value = 2638.923
$ws.Range("K122").Value = 7916.768999999999
$ws.Range("M122").Value = -5466.768999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 23989.3
$ws.Range("I38").Value = 24223.5
$ws.Range("J38").Value = 23833.166
$ws.Range("K38").Value = 24223.5
$ws.Range("L38").Value = 23833.166
$ws.Range("M38").Value = -23807.5
$ws.Range("N38").Value = -24665.166
$ws.Range("H86").Value = 372603.78
$ws.Range("J86").Value = 3676.75
$ws.Range("L86").Value = 3676.75
$ws.Range("N86").Value = -5922.75
$ws.Range("H89").Value = 372603.78
$ws.Range("J89").Value = 3676.75
$ws.Range("L89").Value = 18383.75
$ws.Range("N89").Value = -29615.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7098.0835
$ws.Range("I58").Value = 2027.8334
$ws.Range("K58").Value = 2027.8334
$ws.Range("M58").Value = -1824.8334
$ws.Range("H99").Value = 31017.418
$ws.Range("I99").Value = 71032.336
$ws.Range("K99").Value = 71032.336
$ws.Range("M99").Value = -69534.336
$ws.Range("H126").Value = 31017.418
$ws.Range("I126").Value = 71032.336
$ws.Range("K126").Value = 213097.008
$ws.Range("M126").Value = -210627.008
$ws.Range("H132").Value = 27382.77
$ws.Range("I132").Value = 20304.066
$ws.Range("K132").Value = 60912.198
$ws.Range("M132").Value = -58382.198
$ws.Range("H134").Value = 10618.844
$ws.Range("I134").Value = 7109.8096
$ws.Range("J134").Value = 17317.908
$ws.Range("K134").Value = 21329.4288
$ws.Range("L134").Value = 51953.724
$ws.Range("M134").Value = -18794.4288
$ws.Range("N134").Value = -57023.724
$ws.Range("H136").Value = 7098.0835
$ws.Range("I136").Value = 2027.8334
$ws.Range("K136").Value = 6083.5002
$ws.Range("M136").Value = -3533.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3065.7144
$ws.Range("I131").Value = 1361.5
$ws.Range("J131").Value = 3747.4
$ws.Range("K131").Value = 4084.5
$ws.Range("L131").Value = 11242.2
$ws.Range("M131").Value = 955.5
$ws.Range("N131").Value = -21322.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("H5").Value = 18396.1
$ws.Range("I5").Value = 14499.5
$ws.Range("J5").Value = 19370.25
$ws.Range("K5").Value = 14499.5
$ws.Range("L5").Value = 19370.25
$ws.Range("M5").Value = -14387.5
$ws.Range("N5").Value = -19594.25
$ws.Range("H9").Value = 4666.3335
$ws.Range("I9").Value = 1999
$ws.Range("K9").Value = 1999
$ws.Range("M9").Value = -1829
$ws.Range("H10").Value = 6461.6665
$ws.Range("I10").Value = 5754
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 5754
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -5585
$ws.Range("N10").Value = -10338
$ws.Range("H12").Value = 9000
$ws.Range("I12").Value = 8000
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -7860
$ws.Range("N12").Value = -10280
$ws.Range("H13").Value = 17427.475
$ws.Range("I13").Value = 16765.166
$ws.Range("J13").Value = 18562.857
$ws.Range("K13").Value = 16765.166
$ws.Range("L13").Value = 18562.857
$ws.Range("M13").Value = -16626.166
$ws.Range("N13").Value = -18840.857
$ws.Range("H14").Value = 1434286.1
$ws.Range("I14").Value = 2005600.6
$ws.Range("J14").Value = 6000
$ws.Range("K14").Value = 2005600.6
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = -2005432.6
$ws.Range("N14").Value = -6336
$ws.Range("H80").Value = 3859.6
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 3859.6
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H122").Value = 2632.8
$ws.Range("J122").Value = 2580.8333
$ws.Range("L122").Value = 7742.499899999999
$ws.Range("N122").Value = -12642.4999
$ws.Range("H126").Value = 2987.4285
$ws.Range("J126").Value = 2987.4285
$ws.Range("L126").Value = 8962.2855
$ws.Range("N126").Value = -13902.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4405.091
$ws.Range("I122").Value = 4001
$ws.Range("J122").Value = 5112.25
$ws.Range("K122").Value = 12003
$ws.Range("L122").Value = 15336.75
$ws.Range("M122").Value = -9553
$ws.Range("N122").Value = -20236.75
$ws.Range("H136").Value = 22228242
$ws.Range("I136").Value = 27781552
$ws.Range("K136").Value = 83344656
$ws.Range("M136").Value = -83342106

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 14999.5
$ws.Range("I29").Value = 5000
$ws.Range("K29").Value = 5000
$ws.Range("M29").Value = -4710
$ws.Range("H42").Value = 54499.5
$ws.Range("I42").Value = 45399.6
$ws.Range("K42").Value = 45399.6
$ws.Range("M42").Value = -45021.6
$ws.Range("H43").Value = 39998.332
$ws.Range("I43").Value = 39998.332
$ws.Range("K43").Value = 39998.332
$ws.Range("M43").Value = -39849.332
$ws.Range("H132").Value = 22097.277
$ws.Range("I132").Value = 16291.212
$ws.Range("J132").Value = 35783
$ws.Range("K132").Value = 48873.636
$ws.Range("L132").Value = 107349
$ws.Range("M132").Value = -46343.636
$ws.Range("N132").Value = -112409
$ws.Range("H136").Value = 2680.2727
$ws.Range("I136").Value = 2648.3
$ws.Range("K136").Value = 7944.900000000001
$ws.Range("M136").Value = -5394.900000000001

